$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.06327293813228607,0.9886227250099182,0.01683990657329559,0.9982154965400696),
    @(0.0110556585714221,0.9980016350746155,0.01687898300588131,0.9984385371208191),
    @(0.006905577145516872,0.9982947707176208,0.03968919813632965,0.9984385371208191),
    @(0.003651242004707456,0.9989075660705566,0.003612543921917677,0.9992192983627319),
    @(0.001361945876851678,0.9994937777519226,0.003926234319806099,0.9994423389434814),
    @(0.001570164342410862,0.9996269941329956,0.002773643704131246,0.9995538592338562),
    @(0.0005301287164911628,0.9997868537902832,0.002318680752068758,0.9997769594192505),
    @(0.00170036347117275,0.9996802806854248,0.003270498011261225,0.9997769594192505),
    @(0.0001523441897006705,0.9999467134475708,0.002601358108222485,0.9997769594192505),
    @(0.0007136623607948422,0.9998134970664978,0.003124229377135634,0.9997769594192505),
    @(0.0007594649214297533,0.9997069239616394,0.01760998740792274,0.992304265499115),
    @(0.0005252372939139605,0.9999200701713562,0.003199307480826974,0.9997769594192505),
    @(0.0004926143446937203,0.9998134970664978,0.003098492277786136,0.9997769594192505),
    @(0.000819545763079077,0.9997868537902832,0.003536414820700884,0.9997769594192505),
    @(0.0002019653766183183,0.9999467134475708,0.004196121357381344,0.9997769594192505),
    @(0.0005274725263006985,0.9998934268951416,0.002796070417389274,0.9997769594192505),
    @(0.0003469176299404353,0.9999467134475708,0.003274097573012114,0.9997769594192505),
    @(0.00007005669613135979,0.9999733567237854,0.003827110165730119,0.9997769594192505),
    @(0.0001225347368745133,0.9999467134475708,0.003550265450030565,0.9997769594192505),
    @(0.0004326587950345129,0.9998934268951416,0.003471905365586281,0.9997769594192505),
    @(0.001067671691998839,0.999866783618927,0.004016472958028316,0.999665379524231),
    @(0.0006713774055242538,0.999866783618927,0.003246141597628593,0.9997769594192505),
    @(0.0001092476595658809,0.9999733567237854,0.004160674288868904,0.9995538592338562),
    @(0.0001526341948192567,0.9999733567237854,0.00329888123087585,0.9997769594192505),
    @(0.00008568457269575447,0.9999733567237854,0.003834961215034127,0.9997769594192505),
    @(0.0000426693404733669,0.9999733567237854,0.002556911902502179,0.9997769594192505),
    @(0.0001664258743403479,0.9999467134475708,0.00439690425992012,0.9997769594192505),
    @(0.0004350739473011345,0.999866783618927,0.006373071111738682,0.9997769594192505),
    @(0.0001394504506606609,0.9999467134475708,0.006421537604182959,0.9997769594192505),
    @(0.000013008558198635,1,0.006306247785687447,0.9997769594192505),
    @(0.0003387824981473386,0.9999467134475708,0.00570655195042491,0.9997769594192505),
    @(0.0001498584315413609,0.9999733567237854,0.005583497229963541,0.9997769594192505),
    @(0.0003023763129021972,0.9999200701713562,0.003636211855337024,0.9997769594192505),
    @(0.000009151679478236474,1,0.004160190932452679,0.9997769594192505),
    @(0.000006242999461392174,1,0.00543748727068305,0.9997769594192505),
    @(0.0007000568439252675,0.9997868537902832,0.007759205531328917,0.9997769594192505),
    @(0.00005220238017500378,0.9999733567237854,0.008040698245167732,0.9997769594192505),
    @(0.0002016634098254144,0.9999200701713562,0.006658358499407768,0.9997769594192505),
    @(0.00005523686559172347,0.9999733567237854,0.008671056479215622,0.9997769594192505),
    @(0.00001530833651486319,1,0.008674080483615398,0.9997769594192505),
    @(0.0001030903804348782,0.9999733567237854,0.008542068302631378,0.9997769594192505),
    @(0.0001198840327560902,0.9999733567237854,0.007412600796669722,0.9997769594192505),
    @(0.0009276361670345068,0.9999200701713562,0.007634813897311687,0.9997769594192505),
    @(0.0003222703526262194,0.9999200701713562,0.006075161509215832,0.9997769594192505),
    @(0.0003632468869909644,0.9999467134475708,0.006231239065527916,0.9997769594192505),
    @(0.001047316007316113,0.9998134970664978,0.007820194587111473,0.9997769594192505),
    @(0.0003526316431816667,0.9998934268951416,0.007436858955770731,0.9997769594192505),
    @(0.00007197929517133161,0.9999467134475708,0.008635526522994041,0.9997769594192505),
    @(0.00002932405004685279,1,0.008642366155982018,0.9997769594192505),
    @(0.00001847551720857155,1,0.009329378604888916,0.9997769594192505)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}
